# "Moved code windows to better fit slide."
#
# Slide 9 has a "code window" rectangle ("Rectangle 4") that sits a hair
# above its companion "TextBox 7" (both are part of the same code-sample
# callout). Nudge the rectangle down by 330 EMU so its top edge lines up
# exactly with the textbox's top edge (541020 EMU == 42.6 pt).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item("Rectangle 4")

# Shape.Top is expressed in points (1 pt == 12700 EMU). The target offset
# is 541020 EMU == 42.6 pt exactly; nudge slightly past .6 so float
# round-trip truncation still lands on 541020 EMU rather than 541019.
$shp.Top = 42.60005
